# Add Rev D PCBA files
#
# - Adds a new first sheet "TLV3542 Noise Figure" (placeholder with a
#   reference link + the standard Parameter/Value/Units/Note header row).
# - Renames the existing "Sheet1" (AD8314 input-match calcs) to
#   "AD8314 Input Match" and moves it to the second tab position.
# - Tidies a few leftover per-cell format overrides on the AD8314 sheet
#   that are no longer needed.

$wb = $excel.ActiveWorkbook

# --- Locate / reorder sheets -------------------------------------------------
$existing = $wb.Worksheets.Item("Sheet1")

# Worksheets.Add() with a single "before" argument inserts the new sheet
# immediately before the referenced one, i.e. as the new first tab.
$noiseSheet = $wb.Worksheets.Add($existing)

# Re-fetch the original sheet by name - the freshly inserted sheet now
# occupies position 1, so grab "Sheet1" again by name to get a handle on the
# original (now second) tab.
$matchSheet = $wb.Worksheets.Item("Sheet1")

$matchSheet.Name = "AD8314 Input Match"
$noiseSheet.Name = "TLV3542 Noise Figure"

# --- Populate the new "TLV3542 Noise Figure" sheet --------------------------
$noiseSheet.Range("A1").Value = "Noise figure calculation paper: https://www.ti.com/lit/an/slyt094/slyt094.pdf"
$noiseSheet.Range("A2").Value = "Parameter"
$noiseSheet.Range("B2").Value = "Value"
$noiseSheet.Range("C2").Value = "Units"
$noiseSheet.Range("D2").Value = "Note"

# Leave the whole title row selected, mirroring the tab being left on this
# placeholder sheet after adding it.
$excel.Goto($noiseSheet.Range("A1:XFD1"))

# --- Tidy up stray per-cell formatting on the AD8314 sheet -------------------
# Rows 17-22 carried redundant per-cell style overrides (bold/fill flags that
# were no-ops given the cells' actual font); drop them back to plain cells,
# then restore the wrap-text note formatting in column D for that block.
$matchSheet.Rows("17:22").ClearFormats()
$matchSheet.Range("D17:D22").WrapText = $true

# B14 had a stray "apply fill" flag (with no actual fill colour) left over
# from earlier edits; clear it so the cell is back to the plain default style.
$matchSheet.Range("B14").Interior.Pattern = -4142
